$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains its original plain-text formatting so that
# numeric-looking values (e.g. "212.30", "1.00") are not auto-converted into
# real numbers by Excel, which would corrupt the text content / trailing zeros.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.409.91'
$ws.Range("E2").Value = '  +0.59%  '

$ws.Range("D3").Value = '1.607.55'
$ws.Range("E3").Value = '  +1.00%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = '212.30'
$ws.Range("E5").Value = '  -0.10%  '

$ws.Range("E6").Value = '  -0.28%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("E8").Value = '  -0.17%  '

$ws.Range("E9").Value = '  +0.01%  '

$ws.Range("D10").Value = '19.31'
$ws.Range("E10").Value = '  +1.37%  '

$ws.Range("D11").Value = '0.0857'
$ws.Range("E11").Value = '  +0.59%  '

$ws.Range("D12").Value = '1.834.62'
$ws.Range("E12").Value = '  +1.01%  '

$ws.Range("D13").Value = '1.602.14'
$ws.Range("E13").Value = '  +0.53%  '

$ws.Range("E14").Value = '  -0.04%  '

$ws.Range("D15").Value = '0.508'
$ws.Range("E15").Value = '  -0.25%  '

$ws.Range("D16").Value = '63.56'
$ws.Range("E16").Value = '  -0.44%  '

$ws.Range("D17").Value = '233.42'
$ws.Range("E17").Value = '  +8.16%  '

$ws.Range("D18").Value = '26.415.77'
$ws.Range("E18").Value = '  +0.66%  '

$ws.Range("D19").Value = '7.71'
$ws.Range("E19").Value = '  +5.25%  '

$ws.Range("D20").Value = '0.0₃0725'
$ws.Range("E20").Value = '  -0.28%  '

$ws.Range("E21").Value = '  +0.14%  '

$ws.Range("D22").Value = '4.28'
$ws.Range("E22").Value = '  -0.51%  '

$ws.Range("D23").Value = '8.97'
$ws.Range("E23").Value = '  -0.95%  '

$ws.Range("E24").Value = '  +1.16%  '

$ws.Range("D25").Value = '147.23'
$ws.Range("E25").Value = '  +1.98%  '

$ws.Range("E26").Value = '  +0.05%  '

$ws.Range("D27").Value = '6.98'
$ws.Range("E27").Value = '  +0.25%  '

$ws.Range("E28").Value = '  +1.16%  '

$ws.Range("D29").Value = '15.47'
$ws.Range("E29").Value = '  +2.33%  '

$ws.Range("E30").Value = '  +1.15%  '

$ws.Range("E31").Value = '  +0.29%  '

$ws.Range("D32").Value = '1.492.93'
$ws.Range("E32").Value = '  +5.47%  '

$ws.Range("E33").Value = '  +1.31%  '

$ws.Range("E34").Value = '  -0.54%  '

$ws.Range("E35").Value = '  -0.29%  '

$ws.Range("E36").Value = '  +1.06%  '

$ws.Range("D37").Value = '0.564'
$ws.Range("E37").Value = '  -3.29%  '

$ws.Range("E38").Value = '  -0.02%  '

$ws.Range("D39").Value = '0.823'
$ws.Range("E39").Value = '  -0.08%  '

$ws.Range("D40").Value = '5.81'
$ws.Range("E40").Value = '  -0.70%  '

$ws.Range("D42").Value = '2.20'
$ws.Range("E42").Value = '  +2.55%  '

$ws.Range("D43").Value = '0.934'
$ws.Range("E43").Value = '  -4.24%  '

$ws.Range("D44").Value = '1.746.83'
$ws.Range("E44").Value = '  +1.02%  '

$ws.Range("D45").Value = '0.761'
$ws.Range("E45").Value = '  -0.51%  '

$ws.Range("D46").Value = '60.96'
$ws.Range("E46").Value = '  +0.01%  '

$ws.Range("D47").Value = '89.23'
$ws.Range("E47").Value = '  +3.51%  '

$ws.Range("E48").Value = '  +0.70%  '

$ws.Range("D49").Value = '0.0502'
$ws.Range("E49").Value = '  +0.11%  '

$ws.Range("D50").Value = '0.0963'
$ws.Range("E50").Value = '  +0.97%  '

$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").Value = '0.999'
$ws.Range("E51").Value = '  +0.04%  '
